# Update boards in measurements file
# Target: "Boards" worksheet, columns C (Diam,mm) and D (C_measured,pF)
# for rows 2-27 get new measured values; the thin-line cell borders
# around A2:D27 are removed (kept vertical-center / wrap-text alignment
# and the existing number format on C:D); "Boards" becomes the active
# sheet/tab with A14:D14 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boards")

# New measured values: row -> (C = Diam,mm ; D = C_measured,pF)
$newData = @(
    @(2,  75500, 971),
    @(3,  75200, 9810),
    @(4,  75000, 9740),
    @(5,  75800, 9770),
    @(6,  75800, 9750),
    @(7,  74900, 9760),
    @(8,  74800, 1000),
    @(9,  75100, 990),
    @(10, 75100, 1010),
    @(11, 75000, 980),
    @(12, 74900, 990),
    @(13, 75500, 990),
    @(14, 74800, 990),
    @(15, 37600, 971),
    @(16, 37600, 9810),
    @(17, 37600, 9740),
    @(18, 37500, 9770),
    @(19, 37500, 9750),
    @(20, 37600, 9760),
    @(21, 37500, 1000),
    @(22, 37500, 990),
    @(23, 37600, 1010),
    @(24, 37500, 980),
    @(25, 37400, 990),
    @(26, 37800, 990),
    @(27, 37400, 990)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 3).Value2 = $entry[1]
    $ws.Cells.Item($r, 4).Value2 = $entry[2]
}

# Strip the thin border that used to surround the data block (A2:D27);
# the alignment (vertical-center, wrap) and C:D number format stay.
$ws.Range("A2:D27").Borders.LineStyle = -4142

# Make "Boards" the active sheet, with A14:D14 selected (active cell A14).
$ws.Activate()
$ws.Range("A14:D14").Select()
